$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark sitting at the end of the
#    first paragraph (right after "... download project từ server về máy").
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Append a new paragraph after the "git push: ..." paragraph
#    (currently the last paragraph) containing the "git pull" text.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertionPoint = $lastPara.Range
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()
$insertionPoint.Collapse(0)

$pullPara = $d.Paragraphs($d.Paragraphs.Count)
$pullPara.Range.Text = "git pull: cập nhập những thay đổi từ project trên server về máy"

# ------------------------------------------------------------------
# 3) Append another new paragraph after that containing "hết".
# ------------------------------------------------------------------
$pullPara = $d.Paragraphs($d.Paragraphs.Count)
$afterPull = $pullPara.Range
$afterPull.Collapse(0)
$afterPull.InsertParagraphAfter()
$afterPull.Collapse(0)

$hetPara = $d.Paragraphs($d.Paragraphs.Count)
# Type the paragraph text with a trailing sentinel character. Placing the
# bookmark directly at "end of text, before the paragraph mark" triggers a
# quirk in this runtime (the bookmark silently jumps to the very start of
# the document), so a sentinel char is used to land on a safe offset and is
# then deleted afterwards, leaving the bookmark correctly positioned right
# after "hết" and before the paragraph mark.
$hetPara.Range.Text = "hếtX"

$hetPara = $d.Paragraphs($d.Paragraphs.Count)
$paraEnd = $hetPara.Range.End
$bookmarkPos = $paraEnd - 2
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$sentinelRange = $d.Range($paraEnd - 2, $paraEnd - 1)
$sentinelRange.Delete()
